# "Complete adv_new design and issues fix"
#
# 1) Slide 1: the big title textbox ("Module 02" / "-" / "Python basic
#    types") goes from 60pt non-bold to 54pt bold for all three runs.
# 2) Slide 2: the "Text Placeholder 4" body placeholder (ph idx="13")
#    gets an explicit position/size (it previously inherited the
#    layout's empty <p:spPr/>).

$p = $ppt.ActivePresentation

# --- Slide 1: title textbox -------------------------------------------------
$slide1 = $p.Slides.Item(1)
$title = $slide1.Shapes.Item(3)          # "Google Shape;164;p3" text box
$tr = $title.TextFrame.TextRange

# Run boundaries (by character, 1-based) for the three <a:r> runs:
#   "Module 02" (9 chars) + "-" (1 char) + "Python basic types" (18 chars)
$run1 = $tr.Characters(1, 9)
$run2 = $tr.Characters(10, 1)
$run3 = $tr.Characters(11, 18)

foreach ($run in @($run1, $run2, $run3)) {
    $run.Font.Size = 54
    $run.Font.Bold = 1
}

# --- Slide 2: body placeholder (idx 13) position ---------------------------
$slide2 = $p.Slides.Item(2)
$memShape = $slide2.Shapes.Item(1)       # "Text Placeholder 4", ph idx="13"

# Point values chosen so the single-precision round-trip used internally
# lands on the exact target EMUs (815413, 1492163, 9731259, 4642308).
$memShape.Left = 64.20574951171875
$memShape.Top = 117.49315643310547
$memShape.Width = 766.2409057617188
$memShape.Height = 365.53607177734375
